$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.413.78"
$ws.Range("E2").Value = "  +0.67%  "
$ws.Range("D3").Value = "1.638.08"
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("E5").Value = "  +0.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "305.12"
$ws.Range("E6").Value = "  +0.62%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3723"
$ws.Range("E7").Value = "  -1.43%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "52.02"
$ws.Range("E8").Value = "  +0.59%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3621"
$ws.Range("E9").Value = "  -0.39%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.248"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08112"
$ws.Range("E11").Value = "  -0.29%  "
$ws.Range("E12").Value = "  +0.18%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.81"
$ws.Range("E13").Value = "  -0.47%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.594"
$ws.Range("E14").Value = "  -0.16%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001270"
$ws.Range("E15").Value = "  +1.81%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.275"
$ws.Range("E16").Value = "  -1.97%  "
$ws.Range("D17").Value = "1.629.51"
$ws.Range("E17").Value = "  +1.57%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "94.36"
$ws.Range("E18").Value = "  +0.50%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06871"
$ws.Range("E19").Value = "  -0.25%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.10"
$ws.Range("E20").Value = "  -0.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.510"
$ws.Range("E21").Value = "  -0.58%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  +0.09%  "
$ws.Range("D23").Value = "23.413.78"
$ws.Range("E23").Value = "  +0.65%  "
$ws.Range("E24").Value = "  -1.88%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.414"
$ws.Range("E25").Value = "  +0.91%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.024"
$ws.Range("E26").Value = "  +0.46%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.20"
$ws.Range("E27").Value = "  -0.26%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "151.49"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.271"
$ws.Range("E29").Value = "  +0.30%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "135.70"
$ws.Range("E30").Value = "  +1.15%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.284"
$ws.Range("E31").Value = "  -3.85%  "
$ws.Range("D32").Value = "1.808.08"
$ws.Range("E32").Value = "  +1.57%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.755"
$ws.Range("E33").Value = "  -0.36%  "
$ws.Range("E34").Value = "  -1.52%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02834"
$ws.Range("E35").Value = "  +3.66%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.31"
$ws.Range("E36").Value = "  +0.30%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2522"
$ws.Range("E37").Value = "  -0.46%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.07197"
$ws.Range("E38").Value = "  -4.39%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.08771"
$ws.Range("E39").Value = "  -0.44%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.038"
$ws.Range("E40").Value = "  -0.87%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.371"
$ws.Range("E41").Value = "  -0.15%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7030"
$ws.Range("E42").Value = "  -1.29%  "
$ws.Range("E43").Value = "  -0.98%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.09"
$ws.Range("E44").Value = "  +2.99%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6495"
$ws.Range("E45").Value = "  -1.03%  "
$ws.Range("E46").Value = "  +0.12%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.9997"
$ws.Range("E47").Value = "  +0.17%  "
$ws.Range("E48").Value = "  -0.20%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07962"
$ws.Range("E49").Value = "  -0.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "128.65"
$ws.Range("E50").Value = "  -2.81%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.196"
$ws.Range("E51").Value = "  -0.78%  "
